$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.957.63'
$ws.Range("E2").Value = '  +2.19%  '
$ws.Range("D3").Value = '2.584.95'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'521.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").Value = "'138.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.17%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '2.594.00'
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").Value = "'6.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("E13").Value = '  +3.29%  '
$ws.Range("D14").Value = '3.039.13'
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").Value = '58.855.90'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = "'20.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("D17").Value = '2.608.90'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = "'337.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = "'4.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").Value = "'10.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").Value = "'6.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.41%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = "'65.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.21%  '
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").Value = "'7.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '0.0₃0723'
$ws.Range("E30").Value = '  -2.22%  '
$ws.Range("D31").Value = "'5.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.76%  '
$ws.Range("D32").Value = "'1.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("D33").Value = "'18.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("D34").Value = "'149.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("D35").Value = "'3.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.55%  '
$ws.Range("D36").Value = "'1.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").Value = "'36.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("D38").Value = "'1.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("D39").Value = "'0.825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").Value = "'0.805"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.76%  '
$ws.Range("D41").Value = "'3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = "'270.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("D44").Value = "'10.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("D45").Value = "'0.0953"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("D46").Value = "'0.587"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("D47").Value = "'0.0516"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").Value = "'18.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").Value = '1.961.52'
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").Value = "'4.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.16%  '
